$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last status check" timestamp shown in F1
$ws.Range("F1").Value = "Last status check on: 16.02.2022 06:15"

# Update the Makro row (row 5): new price pushes old price into "Old Cena",
# delta becomes a signed text value, and the "Old Datum" column now stores
# a literal timestamp string instead of a date serial value.
$ws.Range("B5").Value = 36.5
$ws.Range("C5").Value = 36.1

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "+0.4"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2022-02-16 06:15:19"
$ws.Range("E5").Style = "Normal"
